$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3: "Email" -> "EmailID"
$ws.Range("B3").Value = "EmailID"

# E21: new note cell "need to change to city"
$ws.Range("E21").Value = "need to change to city"
$ws.Range("E21").Font.Bold = $true
$ws.Range("E21").Interior.Color = 255

# B30 / D30: cleared (values removed)
$ws.Range("B30").Value = ""
$ws.Range("D30").Value = ""

# Sheet view: scroll/selection change
$ws.Range("D30").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
